$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.178.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.762.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.761.16"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.63%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +1.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("E13").Value = "  -0.41%  "

$ws.Range("E14").Value = "  +1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.392.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.761.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.215.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.51%  "

$ws.Range("E18").Value = "  +1.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.39%  "

$ws.Range("E20").Value = "  -1.41%  "

$ws.Range("E21").Value = "  +18.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "494.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.729"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.81%  "

$ws.Range("E24").Value = "  +7.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "

$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("E27").Value = "  +0.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.51%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.00%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.908.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.701.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.84%  "

$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.76%  "

$ws.Range("E39").Value = "  +1.25%  "

$ws.Range("E40").Value = "  +2.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.327"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.04%  "

$ws.Range("E42").Value = "  +5.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "430.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.74%  "

$ws.Range("E45").Value = "  +0.48%  "

$ws.Range("E46").Value = "  +1.16%  "

$ws.Range("E48").Value = "  -0.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.797.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.88%  "

$ws.Range("E51").Value = "  +0.64%  "
